$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# 1) Merge the split runs in the "Glasgow Haskell Compiler ... binary." /
#    "main benefit of a compiler ..." paragraph into single runs, while
#    leaving the lone middle space run on its own.
# -------------------------------------------------------------------------

$sentence1 = "The Glasgow Haskell Compiler is the main compiler for Haskell. The job of the compiler is to transform human-readable source code into machine-readable binary."

$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("The Glasgow Haskell Compiler is the main compiler for Haskell.")) {
        $para = $p
        break
    }
}

$paraStart = $para.Range.Start
$rngSentence1 = $d.Range($paraStart, $paraStart + $sentence1.Length)
$rngSentence1.Find.Execute($sentence1, $true, $false, $false, $false, $false, $true, 1, $false, $sentence1, 2) | Out-Null

# The single space between the two sentences got folded back into the
# surrounding run by the replace above; split it back out into its own
# run (matching the original document) by toggling a character property.
$spacePos = $paraStart + $sentence1.Length
$rngSpace = $d.Range($spacePos, $spacePos + 1)
$rngSpace.Bold = 1
$rngSpace.Bold = 0

# -------------------------------------------------------------------------
# 2) Append the new paragraphs at the end of the document (after the
#    "...it works." paragraph, before the bookmark).
# -------------------------------------------------------------------------

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRng = $lastPara.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

$p1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p1.Range.InsertAfter("The interactive mode of GHC is GHCi. We use the command ghci to start it and :q to exit on terminal.")

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p2.Range.InsertAfter("Note: Prior to version 8 of GHCi, function and variable definitions needed to be prefaced with a let keyword. This is no longer necessary, but many Haskell examples on the web and in older books still include it:")

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r3 = $p3.Range
$r3.InsertAfter("GHCi> let f x = x + x")
$r3.Font.Name = "Courier New"
$r3.Font.NameAscii = "Courier New"
$r3.Font.NameBi = "Courier New"
$r3.Font.Size = 10
$r3.Font.SizeBi = 10

$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r4 = $p4.Range
$r4.InsertAfter("GHCi> f ")
$r4.Font.Name = "Courier New"
$r4.Font.NameAscii = "Courier New"
$r4.Font.NameBi = "Courier New"
$r4.Font.Size = 10
$r4.Font.SizeBi = 10

$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r5 = $p5.Range
$r5.InsertAfter("2")
$r5.Font.Name = "Courier New"
$r5.Font.NameAscii = "Courier New"
$r5.Font.NameBi = "Courier New"
$r5.Font.Size = 10
$r5.Font.SizeBi = 10

$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r6 = $p6.Range
$r6.Collapse(0)
$r6.Font.Name = "Courier New"
$r6.Font.NameAscii = "Courier New"
$r6.Font.NameBi = "Courier New"
$r6.Font.Size = 10
$r6.Font.SizeBi = 10

Write-Host "Edit complete"
